$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "1.000", "0.000007952") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) values for the refreshed crypto data
$ws.Range("D2").Value = "26.554.85"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.841.66"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "259.99"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.5242"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").Value = "0.3176"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "0.06796"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "18.74"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "0.7827"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "0.07773"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "1.830.53"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "87.95"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "5.020"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "0.000007952"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "26.568.02"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "2.072.51"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "4.613"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "5.982"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "9.333"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "142.66"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "2.213"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "1.675"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").Value = "16.90"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "111.97"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "4.175"
$ws.Range("D31").Value = "0.08723"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "4.079"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "0.04891"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "3.092"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "2.231"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "0.4814"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "0.8961"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").Value = "110.31"
$ws.Range("D43").Value = "5.916"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "7.647"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "0.4171"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").Value = "9.006"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "0.05836"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "0.1232"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").Value = "34.89"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "0.8936"
$ws.Range("E51").Value = "  +1.24%  "

# Row 34 and 35 swapped positions: ARBITRUM now ranks above ImmutableX
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.134"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7233"
$ws.Range("E35").Value = "  +3.66%  "
